# Updating the Quiz Bank
# Fill in "Can't Say" (col E) and "None of the Above" (col F) answer options
# for the question rows that were missing them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$CANT_SAY = "Can't Say"
$NONE_ABOVE = "None of the Above"

# Rows where E/F are currently empty-but-styled cells, or brand-new cells that
# should simply inherit whatever style Excel naturally applies (the row's own
# style, or the column's default style for rows with no explicit row style).
$simpleRows = @(13,15,23,24,27,34,39,40,50,51,53,59)
foreach ($r in $simpleRows) {
    $ws.Range("E$r").Value = $CANT_SAY
    $ws.Range("F$r").Value = $NONE_ABOVE
}

# Row 47: plain fill-in, keeps its own row style (s=15) naturally.
$ws.Range("E47").Value = $CANT_SAY
$ws.Range("F47").Value = $NONE_ABOVE

# Row 46: special case - the new E46/F46 cells end up carrying row 47's style
# (s=15) rather than row 46's own style (s=9), matching the source edit.
$ws.Range("E46").Value = $CANT_SAY
$ws.Range("F46").Value = $NONE_ABOVE
$ws.Range("E47:F47").Copy() | Out-Null
$ws.Range("E46:F46").PasteSpecial(-4122) | Out-Null

# Rows 68, 70, 71, 73, 74, 76: brand new E/F cells that end up with *no*
# explicit style (s attribute absent) in the target file.
$noStyleNewRows = @(68,70,71,73,74,76)
foreach ($r in $noStyleNewRows) {
    $ws.Range("E$r").Value = $CANT_SAY
    $ws.Range("F$r").Value = $NONE_ABOVE
    $rng = "E" + $r + ":F" + $r
    $ws.Range($rng).Style = "Normal"
}

# Rows 69, 72, 75, 77, 78: E/F already had values - only their style
# attribute is stripped (values stay the same).
$stripStyleRows = @(69,72,75,77,78)
foreach ($r in $stripStyleRows) {
    $rng = "E" + $r + ":F" + $r
    $ws.Range($rng).Style = "Normal"
}

$excel.CutCopyMode = 0

# Update the active selection to match the saved workbook state.
$ws.Range("B15").Select() | Out-Null
